$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the three runs that make up
#   "Create " + "__init__.py" + ", train_pipeline.py, predict_pipeline.py
#   files under pipeline folder in vs code."
# into a single run with the same text/formatting. A Find/Replace over the
# whole (already-contiguous) visible text naturally collapses the matched
# runs into one run carrying the shared formatting.
# ---------------------------------------------------------------------------
$target = "Create __init__.py, train_pipeline.py, predict_pipeline.py files under pipeline folder in vs code."
$rng = $d.Content
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $target, 2)

# ---------------------------------------------------------------------------
# Change 2: append a new run ", logger.py" right after the existing
# "Write code in exception.py" run (same paragraph, two runs), then add a
# series of brand new list paragraphs after it.
# ---------------------------------------------------------------------------
$flatOpcHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$flatOpcFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rPr = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

$rng = $d.Content
$found = $rng.Find.Execute("Write code in exception.py", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$runXml = '<w:p><w:r>' + $rPr + '<w:t>, logger.py</w:t></w:r></w:p>'
$rng.InsertXML($flatOpcHeader + $runXml + $flatOpcFooter)

# Re-find the anchor paragraph (its text now ends in ", logger.py") so we can
# append the new paragraphs right after it, in order.
$anchor = $d.Content
$found = $anchor.Find.Execute("Write code in exception.py, logger.py", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$para = $anchor.Paragraphs(1)
$newPara = $para.Next()

function Set-ListParagraph([object]$paragraph, [string]$text, [bool]$numbered, [bool]$preserveSpace) {
    $spaceAttr = ""
    if ($preserveSpace) { $spaceAttr = ' xml:space="preserve"' }
    $pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/>'
    if ($numbered) {
        $pPr += '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>'
    }
    $pPr += $rPr + '</w:pPr>'
    $run = '<w:r>' + $rPr + '<w:t' + $spaceAttr + '>' + $text + '</w:t></w:r>'
    $xml = $flatOpcHeader + '<w:p>' + $pPr + $run + '</w:p>' + $flatOpcFooter
    $wholeRange = $d.Range($paragraph.Range.Start, $paragraph.Range.End)
    $wholeRange.InsertXML($xml)
}

Set-ListParagraph $newPara "In vs code terminal execute the following commands" $true $false
$newPara = $newPara.Next()

Set-ListParagraph $newPara "git status" $false $false
$newPara = $newPara.Next()

Set-ListParagraph $newPara "git add ." $false $false
$newPara = $newPara.Next()

# "git push -u origin main" is split across two runs: "g" + "it push -u origin main"
$pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/>' + $rPr + '</w:pPr>'
$run1 = '<w:r>' + $rPr + '<w:t>g</w:t></w:r>'
$run2 = '<w:r>' + $rPr + '<w:t>it push -u origin main</w:t></w:r>'
$xml = $flatOpcHeader + '<w:p>' + $pPr + $run1 + $run2 + '</w:p>' + $flatOpcFooter
$wholeRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$wholeRange.InsertXML($xml)
$newPara = $newPara.Next()

Set-ListParagraph $newPara "Create a folder notebook in the file explorer where your project is stored" $true $false
$newPara = $newPara.Next()

Set-ListParagraph $newPara "Create a data folder under notebook folder. " $true $true
$newPara = $newPara.Next()

Set-ListParagraph $newPara "Download the dataset file into the data folder." $true $false
$newPara = $newPara.Next()

Set-ListParagraph $newPara "Download both ipynb files under notebook folder." $true $false
